# Reorder the comma-separated "Recorded By" values in column G:
# move the last item in the list to the front of the list.
# Cells with a single value (no comma) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $last = $parts[$parts.Count - 1]
            $rest = $parts[0..($parts.Count - 2)]
            $newParts = @($last) + $rest
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
